$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.989.36'
$ws.Range("E2").Value = '  +0.63%  '

$ws.Range("D3").Value = '1.641.16'
$ws.Range("E3").Value = '  +0.42%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.93'
$ws.Range("E5").Value = '  +0.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5106'
$ws.Range("E6").Value = '  +1.96%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2561'
$ws.Range("E8").Value = '  +0.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06369'
$ws.Range("E9").Value = '  +0.77%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.58'
$ws.Range("E10").Value = '  +1.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07771'
$ws.Range("E11").Value = '  +0.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.295'
$ws.Range("E12").Value = '  +1.74%  '

$ws.Range("D13").Value = '1.641.66'
$ws.Range("E13").Value = '  +0.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5449'
$ws.Range("E14").Value = '  +1.22%  '

$ws.Range("D15").Value = '0.0₅7765'
$ws.Range("E15").Value = '  -0.60%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.33'
$ws.Range("E16").Value = '  +0.40%  '

$ws.Range("D17").Value = '25.986.02'
$ws.Range("E17").Value = '  +0.47%  '

$ws.Range("E18").Value = '  -0.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '197.32'
$ws.Range("E19").Value = '  +1.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.439'
$ws.Range("E20").Value = '  +2.47%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.942'
$ws.Range("E21").Value = '  +1.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.056'
$ws.Range("E22").Value = '  +2.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.004'
$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.873'
$ws.Range("E24").Value = '  -0.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.23'
$ws.Range("E25").Value = '  +1.21%  '

$ws.Range("E26").Value = '  +6.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.849'
$ws.Range("E27").Value = '  +1.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.65'
$ws.Range("E28").Value = '  +0.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.236'
$ws.Range("E29").Value = '  +0.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.04941'
$ws.Range("E30").Value = '  +2.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.250'
$ws.Range("E31").Value = '  +0.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.191'
$ws.Range("E32").Value = '  +1.35%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.527'
$ws.Range("E33").Value = '  +0.44%  '

$ws.Range("E34").Value = '  +0.13%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.8941'
$ws.Range("E35").Value = '  +1.46%  '

$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '1.149.62'
$ws.Range("E36").Value = '  +2.41%  '

$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.584'
$ws.Range("E37").Value = '  -0.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5434'
$ws.Range("E38").Value = '  -0.70%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01554'
$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("E40").Value = '  +0.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.525'
$ws.Range("E41").Value = '  -1.08%  '

$ws.Range("E42").Value = '  +6.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8135'
$ws.Range("E43").Value = '  +1.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.99'
$ws.Range("E44").Value = '  -0.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.446'
$ws.Range("E45").Value = '  -3.65%  '

$ws.Range("D46").Value = '1.771.43'
$ws.Range("E46").Value = '  -0.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4522'
$ws.Range("E47").Value = '  +0.12%  '

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.85'
$ws.Range("E48").Value = '  +1.27%  '

$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9987'
$ws.Range("E49").Value = '  -0.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05052'
$ws.Range("E50").Value = '  +0.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("E51").Value = '  -0.44%  '

$wb.Save()
Write-Host "Applied cryptos update"